$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.9
$ws.Range("H2").Value = 2.78
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 2.62
$ws.Range("P2").Value = 1.34
$ws.Range("Q2").Value = 3.3

# Row 3
$ws.Range("Q3").Value = 1.33

# Row 4
$ws.Range("G4").Value = 3.65
$ws.Range("I4").Value = 2.22
$ws.Range("Q4").Value = 1.47
$ws.Range("S4").Value = 2.16
$ws.Range("X4").Value = 32
$ws.Range("Z4").Value = 970
$ws.Range("AA4").Value = 27
$ws.Range("AB4").Value = 22
$ws.Range("AC4").Value = 970
$ws.Range("AD4").Value = 970
$ws.Range("AE4").Value = 970
$ws.Range("AF4").Value = 30
$ws.Range("AI4").Value = 26
$ws.Range("AJ4").Value = 60
$ws.Range("AK4").Value = 32
$ws.Range("AL4").Value = 34
$ws.Range("AM4").Value = 50
$ws.Range("AN4").Value = 970
$ws.Range("AO4").Value = 970

# Row 5
$ws.Range("F5").Value = 3.6
$ws.Range("I5").Value = 2.22

# Row 6
$ws.Range("F6").Value = 1.74
$ws.Range("G6").Value = 2.16
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 8.800000000000001
$ws.Range("J6").Value = 2.84
$ws.Range("K6").Value = 4.6
$ws.Range("Q6").Value = 2.64

# Row 8
$ws.Range("F8").Value = 6.4
$ws.Range("X8").Value = 23
$ws.Range("AB8").Value = 26
$ws.Range("AC8").Value = 12.5
$ws.Range("AE8").Value = 16.5
$ws.Range("AF8").Value = 70
$ws.Range("AG8").Value = 28
$ws.Range("AI8").Value = 36

# Row 9
$ws.Range("P9").Value = 1.99
$ws.Range("Q9").Value = 1.89

# Row 10
$ws.Range("G10").Value = 2.94
$ws.Range("P10").Value = 2.22

# Row 11
$ws.Range("G11").Value = 5.4
$ws.Range("H11").Value = 1.84
$ws.Range("J11").Value = 3.6
$ws.Range("K11").Value = 3.8
$ws.Range("P11").Value = 1.81

# Row 12
$ws.Range("F12").Value = 1.82
$ws.Range("G12").Value = 1.84
$ws.Range("H12").Value = 4.7
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 3.95
$ws.Range("K12").Value = 4.1
$ws.Range("P12").Value = 2.24
$ws.Range("Q12").Value = 1.7

# Row 13
$ws.Range("S13").Value = 2.82
$ws.Range("X13").Value = 27
$ws.Range("AK13").Value = 16

# Row 14
$ws.Range("F14").Value = 1.4
$ws.Range("G14").Value = 1.41
$ws.Range("H14").Value = 9.6
$ws.Range("I14").Value = 10.5
$ws.Range("J14").Value = 5.1
$ws.Range("K14").Value = 5.3
$ws.Range("P14").Value = 2.18
$ws.Range("Q14").Value = 1.74

# Row 15
$ws.Range("G15").Value = 5.2

# Row 17
$ws.Range("G17").Value = 1.51
$ws.Range("J17").Value = 4.8
$ws.Range("K17").Value = 5.4
$ws.Range("P17").Value = 2.5
$ws.Range("R17").Value = 1.61
$ws.Range("T17").Value = 1.78
$ws.Range("X17").Value = 29
$ws.Range("Y17").Value = 34
$ws.Range("Z17").Value = 75
$ws.Range("AB17").Value = 13
$ws.Range("AD17").Value = 32
$ws.Range("AE17").Value = 110
$ws.Range("AG17").Value = 11
$ws.Range("AH17").Value = 25
$ws.Range("AI17").Value = 90
$ws.Range("AJ17").Value = 15
$ws.Range("AK17").Value = 14.5

# Row 18
$ws.Range("F18").Value = 1.92
$ws.Range("G18").Value = 1.95
$ws.Range("H18").Value = 4
$ws.Range("Q18").Value = 1.79

# Row 19
$ws.Range("P19").Value = 2.08
$ws.Range("R19").Value = 1.42
$ws.Range("S19").Value = 3.1
$ws.Range("X19").Value = 18
$ws.Range("Z19").Value = 50
$ws.Range("AE19").Value = 810
$ws.Range("AF19").Value = 10
$ws.Range("AL19").Value = 40

# Row 20
$ws.Range("G20").Value = 2.4
$ws.Range("H20").Value = 3.2
$ws.Range("J20").Value = 3.6
$ws.Range("P20").Value = 2.08
$ws.Range("Q20").Value = 1.8

# Row 21
$ws.Range("I21").Value = 11
$ws.Range("J21").Value = 5.5
$ws.Range("K21").Value = 5.7
$ws.Range("U21").Value = 1.9
$ws.Range("X21").Value = 29
$ws.Range("Y21").Value = 40
$ws.Range("AJ21").Value = 11.5
$ws.Range("AL21").Value = 40

# Row 22
$ws.Range("I22").Value = 6
$ws.Range("P22").Value = 2.46

# Row 23
$ws.Range("F23").Value = 1.79

# Row 24
$ws.Range("H24").Value = 3.85
$ws.Range("P24").Value = 2.22
$ws.Range("Q24").Value = 1.75

# Row 25
$ws.Range("F25").Value = 1.33
$ws.Range("G25").Value = 1.37
$ws.Range("J25").Value = 5.8
$ws.Range("P25").Value = 2.68
$ws.Range("Q25").Value = 1.51

# Row 26
$ws.Range("G26").Value = 2.92
$ws.Range("H26").Value = 3.25
$ws.Range("I26").Value = 3.6
$ws.Range("J26").Value = 2.8
$ws.Range("K26").Value = 3.05
$ws.Range("Q26").Value = 2.88

# Row 27
$ws.Range("G27").Value = 2.72
$ws.Range("H27").Value = 3.35
$ws.Range("I27").Value = 4.1

# Row 29
$ws.Range("G29").Value = 2.18
$ws.Range("I29").Value = 5.2
$ws.Range("J29").Value = 2.96
$ws.Range("K29").Value = 3.2
$ws.Range("Q29").Value = 2.78

# Row 30
$ws.Range("F30").Value = 2.1
$ws.Range("H30").Value = 4.5
$ws.Range("I30").Value = 5
$ws.Range("Q30").Value = 2.82

# Row 31
$ws.Range("F31").Value = 1.85
$ws.Range("G31").Value = 1.97
$ws.Range("H31").Value = 4.4
$ws.Range("I31").Value = 5
$ws.Range("J31").Value = 3.6
$ws.Range("K31").Value = 4
$ws.Range("M31").Value = 1.07
$ws.Range("N31").Value = 3.65
$ws.Range("O31").Value = 1.31
$ws.Range("P31").Value = 1.9
$ws.Range("Q31").Value = 1.95
$ws.Range("R31").Value = 1.34
$ws.Range("S31").Value = 3.3
$ws.Range("T31").Value = 1.79
$ws.Range("U31").Value = 2.04
$ws.Range("W31").Value = 2.02
$ws.Range("X31").Value = 15.5
$ws.Range("Y31").Value = 17
$ws.Range("Z31").Value = 38
$ws.Range("AA31").Value = 130
$ws.Range("AB31").Value = 9.199999999999999
$ws.Range("AC31").Value = 9
$ws.Range("AD31").Value = 19.5
$ws.Range("AE31").Value = 65
$ws.Range("AF31").Value = 12.5
$ws.Range("AG31").Value = 11
$ws.Range("AH31").Value = 20
$ws.Range("AI31").Value = 70
$ws.Range("AJ31").Value = 23
$ws.Range("AK31").Value = 22
$ws.Range("AL31").Value = 40
$ws.Range("AM31").Value = 130
$ws.Range("AN31").Value = 14
$ws.Range("AO31").Value = 75

# Row 32
$ws.Range("H32").Value = 3.8
$ws.Range("I32").Value = 6.2
$ws.Range("J32").Value = 2.98
$ws.Range("P32").Value = 1.56
$ws.Range("Q32").Value = 2.24
